$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, G across rows 2-7 (F is unchanged)
$data = @{
    2 = @{ B = 0.2881169905109251;    C = 0.3048912486333797;   D = 3.223369029078222;   E = 0.5333859586016987;  G = 4.349763226824225 }
    3 = @{ B = 0.00009552326474482342; C = 0.002658071450198252; D = 0.1496068669990043;  E = 0.5333859586016987;  G = 0.685746420315646 }
    4 = @{ B = 3.272327238179451;     C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    5 = @{ B = 3.272327238179451;     C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    6 = @{ B = 0.2881169905109251;    C = 109.9114832445916;    D = 0.7210945179870265;  E = 13.86384647080068;   G = 124.7845412238902 }
    7 = @{ B = 3.272327238179451;     C = 1.626987699542094;    D = 3.223369029078222;   E = 13.86384647080068;   G = 21.98653043760045 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("G$row").Value = $rowVals.G
}
